$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Build the new (descending) list of period codes: 2507, 2506, ..., 1607
#    (the original sheet listed them ascending 1607..2506; the edit adds a
#    new period 2507 and reverses the whole list to descending order).
# ---------------------------------------------------------------------------
$periods = New-Object System.Collections.ArrayList
for ($yy = 25; $yy -ge 16; $yy--) {
    $startM = 12
    $endM = 1
    if ($yy -eq 25) { $startM = 7 }   # most recent period is 2507
    if ($yy -eq 16) { $endM = 7 }     # oldest period is 1607
    for ($m = $startM; $m -ge $endM; $m--) {
        $code = "{0:D2}{1:D2}" -f $yy, $m
        [void]$periods.Add($code)
    }
}
# $periods.Count is 109 (108 originally + the new "2507")

# ---------------------------------------------------------------------------
# 2) Make room for the extra data row: insert a whole row at 124 so the
#    trailing signature rows (old 128/129) shift down to 129/130, matching
#    the new dimension B2:J130.
# ---------------------------------------------------------------------------
$ws.Rows.Item(124).Insert()

# Copy the special "last row" formatting (row 123) onto the new row 124 ...
$ws.Range("B123:J123").Copy()
$ws.Range("B124:J124").PasteSpecial(-4122)  # xlPasteFormats

# ... then restore row 123 back to the regular data-row formatting (copied
# from row 122, the row right above it).
$ws.Range("B122:J122").Copy()
$ws.Range("B123:J123").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Fill rows 16-124 with the worker/period detail rows, newest period
#    first (row 16 = 2507) down to the oldest (row 124 = 1607).
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $periods.Count; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "45504080"
    $ws.Range("D$r").Value = "NELYS RAMONA RODRIGUEZ BLANCO"
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = 32000
    $ws.Range("G$r").Value = 800000
}

# ---------------------------------------------------------------------------
# 4) Update the summary figures: one more period in arrears (109 instead of
#    108) and the total amount owed increases by the new period's 32000.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 3488000
$ws.Range("F13").Value = 109
